# Updates the cryptocurrency price/volume/hour table (rows 2-51) to
# reflect the refreshed GitHub Actions scrape. Only the cells that
# actually changed are touched; each target value is text (prices,
# percentages and the "Hora" counter are stored as strings in the
# sheet), so NumberFormat is forced to "@" (Text) before assignment
# to stop Excel from re-interpreting "302.84", "-1.45%" or "3" as
# numeric/percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "302.84"
    "E2" = "-1.45%"
    "G2" = "3"
    "D3" = "35.63"
    "E3" = "-1.61%"
    "G3" = "3"
    "D4" = "5.037"
    "E4" = "-1.99%"
    "G4" = "3"
    "D5" = "0.07890"
    "E5" = "-3.10%"
    "G5" = "3"
    "D6" = "1.860"
    "E6" = "-3.73%"
    "G6" = "3"
    "D7" = "4.110"
    "E7" = "-1.94%"
    "G7" = "3"
    "E8" = "0.18%"
    "G8" = "3"
    "D9" = "0.9223"
    "G9" = "3"
    "D10" = "0.1376"
    "E10" = "-0.66%"
    "G10" = "3"
    "E11" = "-1.68%"
    "G11" = "3"
    "D12" = "0.09132"
    "E12" = "-1.14%"
    "G12" = "3"
    "D13" = "0.03471"
    "E13" = "-2.34%"
    "G13" = "3"
    "D14" = "0.09837"
    "E14" = "0.00%"
    "G14" = "3"
    "D15" = "0.001412"
    "E15" = "-0.07%"
    "G15" = "3"
    "D16" = "0.006220"
    "E16" = "5.37%"
    "G16" = "3"
    "D17" = "3.734"
    "E17" = "4.84%"
    "G17" = "3"
    "E18" = "12.27%"
    "G18" = "3"
    "E19" = "0.05%"
    "G19" = "3"
    "E20" = "3.02%"
    "G20" = "3"
    "D21" = "5.160"
    "E21" = "5.12%"
    "G21" = "3"
    "D22" = "0.2204"
    "E22" = "-11.78%"
    "G22" = "3"
    "D23" = "0.04410"
    "E23" = "-2.70%"
    "G23" = "3"
    "D24" = "0.001236"
    "E24" = "1.87%"
    "G24" = "3"
    "D25" = "0.004623"
    "E25" = "-5.37%"
    "G25" = "3"
    "D26" = "0.0001302"
    "E26" = "4.93%"
    "G26" = "3"
    "E27" = "0.15%"
    "G27" = "3"
    "G28" = "3"
    "G29" = "3"
    "G30" = "3"
    "G31" = "3"
    "G32" = "3"
    "G33" = "3"
    "G34" = "3"
    "G35" = "3"
    "G36" = "3"
    "G37" = "3"
    "G38" = "3"
    "D39" = "0.01931"
    "E39" = "-3.69%"
    "G39" = "3"
    "D40" = "0.05071"
    "E40" = "2.79%"
    "G40" = "3"
    "D41" = "0.007579"
    "E41" = "-1.15%"
    "G41" = "3"
    "D42" = "0.01016"
    "E42" = "-8.63%"
    "G42" = "3"
    "D43" = "0.1344"
    "E43" = "-2.60%"
    "G43" = "3"
    "D44" = "0.002163"
    "E44" = "2.94%"
    "G44" = "3"
    "D45" = "0.009723"
    "E45" = "-8.50%"
    "G45" = "3"
    "D46" = "0.00006173"
    "E46" = "-4.40%"
    "G46" = "3"
    "E47" = "0.12%"
    "G47" = "3"
    "D48" = "65.22"
    "E48" = "0.85%"
    "G48" = "3"
    "E49" = "39.52%"
    "G49" = "3"
    "E50" = "0.12%"
    "G50" = "3"
    "E51" = "0.12%"
    "G51" = "3"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells"
